# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the "d4ea906c-8953-415d-a8ac-ef6a5f9e6338" file, adding
# detailed error messages for the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$zhCnError = "Handback file name: ehbdkfoe.m40 is different with handoff file name: d4ea906c-8953-415d-a8ac-ef6a5f9e6338.6c448c083ad407c0fdba874c56173d408c6ed6e7.zh-cn."
$deDeError = "Handback file name: ehbdkfoe.m40 is different with handoff file name: d4ea906c-8953-415d-a8ac-ef6a5f9e6338.6c448c083ad407c0fdba874c56173d408c6ed6e7.de-de."
$newStatus = "Handback transform failed"

# --- Overview sheet: status column reflects the shared "Status" string ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Note: Excel's COM ColumnWidth setter quantizes to whole pixels and the
# saved OOXML width is recomputed as (round(width*6)+5)/6, so to persist an
# exact width of 40 we must request 39.1666666666667 (=235/6) as input.
$colPWidthForSaved40 = 39.1666666666667

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = $zhCnError
$wsZhCn.Columns.Item(16).ColumnWidth = $colPWidthForSaved40

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = $deDeError
$wsDeDe.Columns.Item(16).ColumnWidth = $colPWidthForSaved40
